# Updates the "Linea 141" schedule workbook with the latest scrape results.
# - Sheet "LP1912": refresh header metadata and append/merge the newest
#   arrivals (rows 38-59).
# - Sheet "LP1912-215": refresh header metadata and append the one new
#   "215"-route arrival that also showed up on the LP1912 sheet.
# - Sheet "6203-6173": refresh the global "last updated" stamp only; no new
#   6203/6173-route arrivals were scraped in this run.

$wb = $excel.ActiveWorkbook

$wsLP1912     = $wb.Worksheets.Item("LP1912")
$wsLP1912_215 = $wb.Worksheets.Item("LP1912-215")
$wsL6203_6173 = $wb.Worksheets.Item("6203-6173")

# ---------------------------------------------------------------------------
# Sheet "LP1912"
# ---------------------------------------------------------------------------
$wsLP1912.Range("A2").Value = "Última actualización: 06:44:40"
$wsLP1912.Range("A3").Value = "Total filas: 54"

$lp1912Rows = @(
  ,@(38, "06:44:40", "07:05", "23_HERNANDEZ", 21, "LP1912")
  ,@(39, "05:18:42", "07:07", "225_GOMEZ", 109, "LP1912")
  ,@(40, "06:44:40", "07:09", "15_ABASTO", 25, "LP1912")
  ,@(41, "05:18:42", "07:11", "215A_EL PATO", 113, "LP1912")
  ,@(42, "05:18:42", "07:15", "11_ETCHEVERRY", 117, "LP1912")
  ,@(43, "06:44:40", "07:16", "16_SANTA ANA", 32, "LP1912")
  ,@(44, "05:53:46", "07:21", "26_HERNANDEZ", 88, "LP1912")
  ,@(45, "06:15:33", "07:23", "10_OLMOS", 68, "LP1912")
  ,@(46, "05:53:46", "07:31", "11_ETCHEVERRY", 98, "LP1912")
  ,@(47, "05:53:46", "07:32", "84_COLONIA URQUIZA-ESC 49", 99, "LP1912")
  ,@(48, "05:53:46", "07:36", "27_EL RETIRO", 103, "LP1912")
  ,@(49, "05:53:46", "07:39", "10_OLMOS", 106, "LP1912")
  ,@(50, "05:53:46", "07:47", "14_ABASTO", 114, "LP1912")
  ,@(51, "05:53:46", "07:51", "215D_EL PATO", 118, "LP1912")
  ,@(52, "06:15:33", "08:07", "16_SANTA ANA", 112, "LP1912")
  ,@(53, "06:15:33", "08:12", "15_ABASTO", 117, "LP1912")
  ,@(54, "06:44:40", "08:21", "26_HERNANDEZ", 97, "LP1912")
  ,@(55, "06:44:40", "08:22", "16_P MOR-SANTA ANA", 98, "LP1912")
  ,@(56, "06:44:40", "08:23", "215B_EL PATO", 99, "LP1912")
  ,@(57, "06:44:40", "08:27", "84_COLONIA URQUIZA-ESC 49", 103, "LP1912")
  ,@(58, "06:44:40", "08:35", "23_HERNANDEZ", 111, "LP1912")
  ,@(59, "06:44:40", "08:42", "81_EL PELIGRO", 118, "LP1912")
)

foreach ($r in $lp1912Rows) {
  $rowNum = $r[0]
  $wsLP1912.Cells.Item($rowNum, 1).Value = $r[1]
  $wsLP1912.Cells.Item($rowNum, 2).Value = $r[2]
  $wsLP1912.Cells.Item($rowNum, 3).Value = $r[3]
  $wsLP1912.Cells.Item($rowNum, 4).Value = $r[4]
  $wsLP1912.Cells.Item($rowNum, 5).Value = $r[5]
}

# ---------------------------------------------------------------------------
# Sheet "LP1912-215"
# ---------------------------------------------------------------------------
$wsLP1912_215.Range("A2").Value = "Última actualización: 06:44:40"
$wsLP1912_215.Range("A3").Value = "Total filas: 13"

$wsLP1912_215.Cells.Item(18, 1).Value = "06:44:40"
$wsLP1912_215.Cells.Item(18, 2).Value = "08:23"
$wsLP1912_215.Cells.Item(18, 3).Value = "215B_EL PATO"
$wsLP1912_215.Cells.Item(18, 4).Value = 99
$wsLP1912_215.Cells.Item(18, 5).Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet "6203-6173"
# ---------------------------------------------------------------------------
$wsL6203_6173.Range("A2").Value = "Última actualización: 06:44:40"
